$wb = $excel.ActiveWorkbook

# --- sheet1 (Worksheets.Item(1)) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 1307
$ws.Range("F7").Value = 7635
$ws.Range("F8").Value = 94
$ws.Range("F11").Value = 8352
$ws.Range("F15").Value = 5647
$ws.Range("F17").Value = 2604
$ws.Range("F20").Value = 342
$ws.Range("F24").Value = 523
$ws.Range("F25").Value = 3463
$ws.Range("F27").Value = 38
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "北京·幻兽帕鲁only"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "北京展览馆 北京展览馆"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "2024.04.04 09:30-04.05 17:00"
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 80
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=82549"
$ws.Range("I29").NumberFormat = "@"
$ws.Range("I29").Value = "//i0.hdslb.com/bfs/openplatform/202403/BbKUlDVR1709866539810.jpeg"
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "2024-04-04"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "北京·第15届IJOY漫展xCGF游戏节"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "北京国家会议中心 北京国家会议中心"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "2024.04.04 09:00-04.05 17:00"
$ws.Range("F30").Value = 2955
$ws.Range("G30").Value = 8.800000000000001
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=81174"
$ws.Range("I30").NumberFormat = "@"
$ws.Range("I30").Value = "//i0.hdslb.com/bfs/openplatform/202401/EJejgoZa1705892035599.jpeg"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "北京·IDOx梦次元动漫游戏嘉年华3rd·配音演员 小N&小敢 专场活动"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "2024.04.05 10:30-04.05 13:45"
$ws.Range("F31").Value = 11
$ws.Range("G31").Value = 268
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=82531"
$ws.Range("I31").NumberFormat = "@"
$ws.Range("I31").Value = "//i1.hdslb.com/bfs/openplatform/202403/S270ineo1709807616493.png"
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "2024-04-05"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "北京·IDOx梦次元动漫游戏嘉年华3rd·配音演员 杨天翔 专场活动"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "北京展览馆 北京展览馆"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "2024.04.05 12:00-04.05 15:00"
$ws.Range("F32").Value = 67
$ws.Range("G32").Value = 258
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "https://show.bilibili.com/platform/detail.html?id=82502"
$ws.Range("I32").NumberFormat = "@"
$ws.Range("I32").Value = "//i0.hdslb.com/bfs/openplatform/202403/pW5dqxbn1709797487963.png"
$ws.Range("F33").Value = 339
$ws.Range("F34").Value = 127
$ws.Range("F35").Value = 302
$ws.Range("F36").Value = 200
$ws.Range("G36").Value = 13.5
$ws.Range("F39").Value = 878
$ws.Range("F40").Value = 1667
$ws.Range("F41").Value = 46
$ws.Range("F43").Value = 14
$ws.Range("F44").Value = 2716
$ws.Range("F46").Value = 2280

# --- sheet2 (Worksheets.Item(2)) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F9").Value = 113

# --- sheet3 (Worksheets.Item(3)) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 265
$ws.Range("F3").Value = 1322

# --- sheet4 (Worksheets.Item(4)) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 265
$ws.Range("F3").Value = 1322
$ws.Range("F5").Value = 1307
$ws.Range("F6").Value = 7635
$ws.Range("F7").Value = 94
$ws.Range("F10").Value = 8356
$ws.Range("F14").Value = 5647
$ws.Range("F16").Value = 2604
$ws.Range("F24").Value = 523
$ws.Range("F26").Value = 3463
$ws.Range("F27").Value = 38
$ws.Range("F29").Value = 2955
$ws.Range("F30").Value = 339
$ws.Range("F31").Value = 127
$ws.Range("F32").Value = 302
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "2024-04-19"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "北京·第22届中国国际模型博览会"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "北京展览馆 北京展览馆"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "2024.04.19 10:00-04.21 17:00"
$ws.Range("F34").Value = 200
$ws.Range("G34").Value = 13.5
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "https://show.bilibili.com/platform/detail.html?id=82425"
$ws.Range("I34").NumberFormat = "@"
$ws.Range("I34").Value = "//i2.hdslb.com/bfs/openplatform/202403/9nkCFSHm1709710888611.jpeg"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "北京·QMQ动漫游戏嘉年华"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "小关路39号 北投购物公园"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "2024.04.20 09:00-04.21 17:00"
$ws.Range("F35").Value = 652
$ws.Range("G35").Value = 63
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = "https://show.bilibili.com/platform/detail.html?id=81982"
$ws.Range("I35").NumberFormat = "@"
$ws.Range("I35").Value = "//i0.hdslb.com/bfs/openplatform/202402/lyPb1fLO1708569465126.jpeg"
$ws.Range("F37").Value = 878
$ws.Range("F39").Value = 1667
$ws.Range("F40").Value = 46
$ws.Range("F42").Value = 14
$ws.Range("F43").Value = 2716
$ws.Range("F46").Value = 2280
$ws.Range("F49").Value = 113
